# This script reproduces the commit "culture collection を MIxS から再度削除　INSDC2017 での確認に基づく":
# the 'culture_collection' column (AM) is removed from the MIGS.eu.built.4.0
# header row (row 15) together with its cell comment. Every column to its
# right (AN..BM) shifts one position to the left (AM..BL), and the header +
# comment that used to live in the now-unused trailing column (BM) are
# dropped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift header values and comment text for columns AM..BL left by one column:
# each column takes over the value/comment that used to belong to the next
# column (e.g. AM, which held 'culture_collection', now holds what AN held,
# i.e. 'dew_point', and so on through BL).
$ws.Range("AM15").Value = 'dew_point'
[void]$ws.Range("AM15").Comment.Text('temperature to which a given parcel of humid air must be cooled, at constant barometric pressure, for water vapor to condense into water.')
$ws.Range("AN15").Value = 'extrachrom_elements'
[void]$ws.Range("AN15").Comment.Text('Plasmids that have significance phenotypic consequence')
$ws.Range("AO15").Value = 'health_state'
[void]$ws.Range("AO15").Comment.Text('Health or disease status of sample at time of collection')
$ws.Range("AP15").Value = 'host'
[void]$ws.Range("AP15").Comment.Text('The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, "Homo sapiens".')
$ws.Range("AQ15").Value = 'host_taxid'
[void]$ws.Range("AQ15").Comment.Text('NCBI taxonomy ID of the host, e.g. 9606')
$ws.Range("AR15").Value = 'indoor_surf'
[void]$ws.Range("AR15").Comment.Text('type of indoor surface')
$ws.Range("AS15").Value = 'isolation_source'
[void]$ws.Range("AS15").Comment.Text('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.')
$ws.Range("AT15").Value = 'locus_tag_prefix'
[void]$ws.Range("AT15").Comment.Text('A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html')
$ws.Range("AU15").Value = 'pathogenicity'
[void]$ws.Range("AU15").Comment.Text('To what is the entity pathogenic')
$ws.Range("AV15").Value = 'ref_biomaterial'
[void]$ws.Range("AV15").Comment.Text('Primary publication or genome report in the form of pubmed ID, DOI or URL')
$ws.Range("AW15").Value = 'samp_collect_device'
[void]$ws.Range("AW15").Comment.Text('Method or device employed for collecting sample')
$ws.Range("AX15").Value = 'samp_mat_process'
[void]$ws.Range("AX15").Comment.Text('Processing applied to the sample during or after isolation')
$ws.Range("AY15").Value = 'samp_size'
[void]$ws.Range("AY15").Comment.Text('Amount or size of sample (volume, mass or area) that was collected')
$ws.Range("AZ15").Value = 'samp_sort_meth'
[void]$ws.Range("AZ15").Comment.Text('method by which samples are sorted')
$ws.Range("BA15").Value = 'samp_vol_we_dna_ext'
[void]$ws.Range("BA15").Comment.Text('volume (mL) or weight (g) of sample processed for DNA extraction')
$ws.Range("BB15").Value = 'source_material_id'
[void]$ws.Range("BB15").Comment.Text('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.')
$ws.Range("BC15").Value = 'specimen_voucher'
[void]$ws.Range("BC15").Comment.Text('Identifier for the physical specimen. Use format: "[<institution-code>:[<collection-code>:]]<specimen_id>", eg, "UAM:Mamm:52179". Intended as a reference to the physical specimen that remains after it was analyzed. If the specimen was destroyed in the process of analysis, electronic images (e-vouchers) are an adequate substitute for a physical voucher specimen. Ideally the specimens will be deposited in a curated museum, herbarium, or frozen tissue collection, but often they will remain in a personal or laboratory collection for some time before they are deposited in a curated collection. There are three forms of specimen_voucher qualifiers. If the text of the qualifier includes one or more colons it is a ''structured voucher''. Structured vouchers include institution-codes (and optional collection-codes) taken from a controlled vocabulary maintained by the INSDC that denotes the museum or herbarium collection where the specimen resides, please visit the INSDC website, http://www.insdc.org/controlled-vocabulary-specimenvoucher-qualifier')
$ws.Range("BD15").Value = 'subspecf_gen_lin'
[void]$ws.Range("BD15").Comment.Text('Information about the genetic distinctness of the lineage (eg., biovar, serovar)')
$ws.Range("BE15").Value = 'substructure_type'
[void]$ws.Range("BE15").Comment.Text('substructure or under building is that largely hidden section of the building which is built off the foundations to the ground floor level')
$ws.Range("BF15").Value = 'surf_air_cont'
[void]$ws.Range("BF15").Comment.Text('contaminant identified on surface')
$ws.Range("BG15").Value = 'surf_humidity'
[void]$ws.Range("BG15").Comment.Text('surfaces: water activity as a function of air and material moisture')
$ws.Range("BH15").Value = 'surf_material'
[void]$ws.Range("BH15").Comment.Text('surface materials at the point of sampling')
$ws.Range("BI15").Value = 'surf_moisture'
[void]$ws.Range("BI15").Comment.Text('water held on a surface')
$ws.Range("BJ15").Value = 'surf_moisture_ph'
[void]$ws.Range("BJ15").Comment.Text('pH measurement of surface')
$ws.Range("BK15").Value = 'surf_temp'
[void]$ws.Range("BK15").Comment.Text('temperature of the surface at the time of sampling')
$ws.Range("BL15").Value = 'trophic_level'
[void]$ws.Range("BL15").Comment.Text('Feeding position in food chain (eg., chemolithotroph)')

# The last column (BM15) now duplicates what BL15 holds; remove its comment
# and shift it out of row 15 entirely so the row once again ends at BL15.
$ws.Range("BM15").Comment.Delete()
$ws.Range("BM15").Delete(-4159)

